# Update countries & provincias Spain
# Applies the daily data refresh: reorders a few tied-rank countries
# (by writing their swapped names/values directly) and updates case
# counts for many rows, plus refreshes the "datos actualizados" timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

$ws.Cells.Item(1, 1).Value = 'Datos actualizados a 18 de Junio de 2020 a las 01:14'
$ws.Cells.Item(4, 2).Value = 2233104
$ws.Cells.Item(4, 3).Value = 24704
$ws.Cells.Item(4, 4).Value = 912372
$ws.Cells.Item(4, 5).Value = 1200802
$ws.Cells.Item(4, 7).Value = 798
$ws.Cells.Item(4, 8).Value = 119930
$ws.Cells.Item(5, 2).Value = 960309
$ws.Cells.Item(5, 3).Value = 31475
$ws.Cells.Item(5, 5).Value = 436280
$ws.Cells.Item(5, 7).Value = 1209
$ws.Cells.Item(5, 8).Value = 46665
$ws.Cells.Item(20, 2).Value = 99853
$ws.Cells.Item(20, 3).Value = 386
$ws.Cells.Item(20, 4).Value = 62017
$ws.Cells.Item(20, 5).Value = 29582
$ws.Cells.Item(26, 1).Value = 'Colombia'
$ws.Cells.Item(26, 2).Value = 57046
$ws.Cells.Item(26, 3).Value = 2115
$ws.Cells.Item(26, 4).Value = 21326
$ws.Cells.Item(26, 5).Value = 33856
$ws.Cells.Item(26, 7).Value = 63
$ws.Cells.Item(26, 8).Value = 1864
$ws.Cells.Item(27, 1).Value = 'Bielorrusia'
$ws.Cells.Item(27, 2).Value = 56032
$ws.Cells.Item(27, 3).Value = 663
$ws.Cells.Item(27, 4).Value = 32735
$ws.Cells.Item(27, 5).Value = 22973
$ws.Cells.Item(27, 7).Value = 6
$ws.Cells.Item(27, 8).Value = 324
$ws.Cells.Item(53, 1).Value = 'Nigeria'
$ws.Cells.Item(53, 2).Value = 17735
$ws.Cells.Item(53, 3).Value = 587
$ws.Cells.Item(53, 4).Value = 5967
$ws.Cells.Item(53, 5).Value = 11299
$ws.Cells.Item(53, 7).Value = 14
$ws.Cells.Item(53, 8).Value = 469
$ws.Cells.Item(54, 1).Value = 'Japon'
$ws.Cells.Item(54, 2).Value = 17628
$ws.Cells.Item(54, 3).Value = 41
$ws.Cells.Item(54, 4).Value = 15850
$ws.Cells.Item(54, 5).Value = 847
$ws.Cells.Item(54, 7).Value = 4
$ws.Cells.Item(54, 8).Value = 931
$ws.Cells.Item(55, 1).Value = 'Austria'
$ws.Cells.Item(55, 2).Value = 17203
$ws.Cells.Item(55, 3).Value = 14
$ws.Cells.Item(55, 4).Value = 16099
$ws.Cells.Item(55, 5).Value = 417
$ws.Cells.Item(55, 7).Value = 6
$ws.Cells.Item(55, 8).Value = 687
$ws.Cells.Item(65, 2).Value = 10162
$ws.Cells.Item(65, 3).Value = 51
$ws.Cells.Item(65, 5).Value = 2430
$ws.Cells.Item(69, 2).Value = 8692
$ws.Cells.Item(69, 3).Value = 32
$ws.Cells.Item(69, 5).Value = 311
$ws.Cells.Item(71, 2).Value = 8020
$ws.Cells.Item(71, 3).Value = 280
$ws.Cells.Item(71, 4).Value = 2966
$ws.Cells.Item(71, 5).Value = 4567
$ws.Cells.Item(71, 7).Value = 10
$ws.Cells.Item(71, 8).Value = 487
$ws.Cells.Item(75, 1).Value = 'Costa de Marfil'
$ws.Cells.Item(75, 2).Value = 6063
$ws.Cells.Item(75, 3).Value = 384
$ws.Cells.Item(75, 4).Value = 2749
$ws.Cells.Item(75, 5).Value = 3266
$ws.Cells.Item(75, 7).Value = 2
$ws.Cells.Item(75, 8).Value = 48
$ws.Cells.Item(76, 1).Value = 'Uzbekistan'
$ws.Cells.Item(76, 2).Value = 5682
$ws.Cells.Item(76, 3).Value = 189
$ws.Cells.Item(76, 4).Value = 4131
$ws.Cells.Item(76, 5).Value = 1532
$ws.Cells.Item(76, 8).Value = 19
$ws.Cells.Item(80, 2).Value = 4668
$ws.Cells.Item(80, 3).Value = 29
$ws.Cells.Item(80, 4).Value = 3364
$ws.Cells.Item(80, 5).Value = 1278
$ws.Cells.Item(91, 1).Value = 'Venezuela'
$ws.Cells.Item(91, 2).Value = 3386
$ws.Cells.Item(91, 3).Value = 236
$ws.Cells.Item(91, 4).Value = 835
$ws.Cells.Item(91, 5).Value = 2523
$ws.Cells.Item(91, 7).Value = 1
$ws.Cells.Item(91, 8).Value = 28
$ws.Cells.Item(92, 1).Value = 'Grecia'
$ws.Cells.Item(92, 2).Value = 3203
$ws.Cells.Item(92, 3).Value = 55
$ws.Cells.Item(92, 4).Value = 1374
$ws.Cells.Item(92, 5).Value = 1642
$ws.Cells.Item(92, 7).Value = 2
$ws.Cells.Item(92, 8).Value = 187
$ws.Cells.Item(98, 2).Value = 2345
$ws.Cells.Item(98, 3).Value = 12
$ws.Cells.Item(98, 4).Value = 2066
$ws.Cells.Item(98, 5).Value = 250
$ws.Cells.Item(130, 2).Value = 899
$ws.Cells.Item(130, 3).Value = 4
$ws.Cells.Item(130, 5).Value = 37
$ws.Cells.Item(135, 4).Value = 810
$ws.Cells.Item(135, 5).Value = 15
$ws.Cells.Item(136, 4).Value = 377
$ws.Cells.Item(136, 5).Value = 408
$ws.Cells.Item(165, 4).Value = 132
$ws.Cells.Item(165, 5).Value = 60
$ws.Cells.Item(169, 4).Value = 102
$ws.Cells.Item(169, 5).Value = 57
$ws.Cells.Item(171, 5).Value = 84
$ws.Cells.Item(171, 7).Value = 1
$ws.Cells.Item(171, 8).Value = 7
$ws.Cells.Item(206, 1).Value = 'Groenlandia'
$ws.Cells.Item(207, 1).Value = 'Islas Malvinas'
$ws.Cells.Item(208, 1).Value = 'Islas Turcas y Caicos'
$ws.Cells.Item(208, 4).Value = 11
$ws.Cells.Item(208, 8).Value = 1
$ws.Cells.Item(209, 1).Value = 'Santa Sede'
$ws.Cells.Item(209, 4).Value = 12
$ws.Cells.Item(209, 8).Value = 0
